$d = $word.ActiveDocument

# The "Yth." (addressee) table row currently reads "Kepala ${satker}".
# Fix the "Jabatan Tujuan Surat" by turning the literal "Kepala " prefix
# into a "${jabatanPimpinan} " placeholder, so the cell text becomes
# "${jabatanPimpinan} ${satker}". Scope the Find to that specific table
# cell so the other, unrelated "Kepala Kantor ..." strings elsewhere in
# the document are left untouched.
$cell = $d.Tables(2).Cell(1, 2)
$cell.Range.Find.Execute("Kepala ", $true, $false, $false, $false,
                          $false, $true, 1, $false,
                          "`${jabatanPimpinan} ", 1)
